$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 138 : compléter l'heure de fin (D138)
$ws.Range("D138").Value = 0.39930555555555558

# Row 139 : nouvelle entrée "Réalisation" / "Commentaires utilisateur"
$ws.Range("A139").Value = 44348
$ws.Range("B139").Value = 5
$ws.Range("C139").Value = 0.40972222222222227
$ws.Range("D139").Value = 0.42430555555555555
$ws.Range("F139").Value = "Réalisation"
$ws.Range("G139").Value = "Commentaires utilisateur"

# Row 140 : nouvelle entrée "Réalisation" / "Note moyenne"
$ws.Range("A140").Value = 44348
$ws.Range("B140").Value = 5
$ws.Range("C140").Value = 0.42430555555555555
$ws.Range("D140").Value = 0.45624999999999999
$ws.Range("F140").Value = "Réalisation"
$ws.Range("G140").Value = "Note moyenne"
$ws.Range("H140").Value = "Affichage de la note moyenne des articles."

# Row 141 : nouvelle entrée "Réalisation" en cours (pas d'heure de fin)
$ws.Range("A141").Value = 44348
$ws.Range("B141").Value = 5
$ws.Range("C141").Value = 0.45624999999999999
$ws.Range("F141").Value = "Réalisation"

# Sélection active déplacée sur G141
$ws.Range("G141").Select()
